$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix 1: comma -> period punctuation in a few "Razon social" entries ---
$nameFixes = @(
    @("E36",  "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"),
    @("E86",  "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"),
    @("E83",  "FERNANDEZ. MARIO HUGO"),
    @("E154", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),
    @("E171", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),
    @("E157", "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN")
)

foreach ($fix in $nameFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# --- Fix 2: re-format "Importe" (column H) amounts from Argentine   ---
# --- "1.234,56" notation to plain decimal "1234.56" notation        ---
$importeFixes = @(
    @(2, "6220.00"),
    @(3, "140000.00"),
    @(4, "8516778.75"),
    @(5, "177000.00"),
    @(6, "756480.00"),
    @(7, "200.00"),
    @(8, "70204.31"),
    @(9, "3872.00"),
    @(10, "5039.65"),
    @(11, "1750.00"),
    @(12, "4043.00"),
    @(13, "7300.00"),
    @(14, "3571.92"),
    @(15, "6600.00"),
    @(16, "525.00"),
    @(17, "12690.00"),
    @(18, "473586.28"),
    @(19, "57176.88"),
    @(20, "95.00"),
    @(21, "15576.00"),
    @(22, "1945.00"),
    @(23, "7934.66"),
    @(24, "642.00"),
    @(25, "15312.00"),
    @(26, "19690.84"),
    @(27, "25308.36"),
    @(28, "2500.00"),
    @(29, "1250.00"),
    @(30, "120.00"),
    @(31, "23577.00"),
    @(32, "1568.37"),
    @(33, "329.42"),
    @(34, "32234.57"),
    @(35, "307.43"),
    @(36, "640.00"),
    @(37, "44.60"),
    @(38, "725.29"),
    @(39, "91.14"),
    @(40, "3989.68"),
    @(41, "2898.26"),
    @(42, "1038.00"),
    @(43, "7099.00"),
    @(44, "35577.80"),
    @(45, "108.00"),
    @(46, "2497.82"),
    @(47, "100.00"),
    @(48, "697.50"),
    @(49, "2050.00"),
    @(50, "26635.14"),
    @(51, "16272.68"),
    @(52, "1477.66"),
    @(53, "985.00"),
    @(54, "53.00"),
    @(55, "251.32"),
    @(56, "239.16"),
    @(57, "2435.60"),
    @(58, "87457.42"),
    @(59, "243.56"),
    @(60, "31735.57"),
    @(61, "2714.55"),
    @(62, "57.96"),
    @(63, "2823.72"),
    @(64, "3434.00"),
    @(65, "1362.00"),
    @(66, "2345.00"),
    @(67, "558.85"),
    @(68, "62389.11"),
    @(69, "1162.00"),
    @(70, "1723.00"),
    @(71, "70.02"),
    @(72, "38300.00"),
    @(73, "125.00"),
    @(74, "136165.55"),
    @(75, "4950.00"),
    @(76, "7115.00"),
    @(77, "1984.80"),
    @(78, "1200.00"),
    @(79, "17071.00"),
    @(80, "79500.00"),
    @(81, "43836.00"),
    @(82, "153738.00"),
    @(83, "40.00"),
    @(84, "17860.00"),
    @(85, "9160.00"),
    @(86, "13680.00"),
    @(87, "611.00"),
    @(88, "2115.00"),
    @(89, "4790.00"),
    @(90, "21.59"),
    @(91, "6100.00"),
    @(92, "391178.31"),
    @(93, "264.00"),
    @(94, "53261.49"),
    @(95, "4235.00"),
    @(96, "4500.00"),
    @(97, "133.48"),
    @(98, "44.00"),
    @(99, "630.00"),
    @(100, "2500.00"),
    @(101, "16734.50"),
    @(102, "54.45"),
    @(103, "3775.00"),
    @(104, "278.98"),
    @(105, "6961.50"),
    @(106, "10830.00"),
    @(107, "2074.00"),
    @(108, "634.10"),
    @(109, "6590.40"),
    @(110, "2997.00"),
    @(111, "288.76"),
    @(112, "1299.00"),
    @(113, "389.90"),
    @(114, "9654.00"),
    @(115, "561.50"),
    @(116, "639.80"),
    @(117, "6130.00"),
    @(118, "5500.00"),
    @(119, "3800.00"),
    @(120, "14100.00"),
    @(121, "7848.00"),
    @(122, "9969.82"),
    @(123, "3260.00"),
    @(124, "1010.00"),
    @(125, "4000.00"),
    @(126, "13000.00"),
    @(127, "2400.00"),
    @(128, "3500.00"),
    @(129, "2800.00"),
    @(130, "13300.00"),
    @(131, "7500.00"),
    @(132, "1759.40"),
    @(133, "2479.00"),
    @(134, "1715.00"),
    @(135, "348.16"),
    @(136, "263102.50"),
    @(137, "1320.00"),
    @(138, "198080.00"),
    @(139, "42797.84"),
    @(140, "2300.00"),
    @(141, "1500.00"),
    @(142, "24478.30"),
    @(143, "1657.50"),
    @(144, "13500.00"),
    @(145, "950.00"),
    @(146, "1250.00"),
    @(147, "346.85"),
    @(148, "3000.00"),
    @(149, "12180.00"),
    @(150, "300.00"),
    @(151, "22950.00"),
    @(152, "937.51"),
    @(153, "6670.00"),
    @(154, "980.00"),
    @(155, "1330.00"),
    @(156, "36.04"),
    @(157, "9060.00"),
    @(158, "2155.00"),
    @(159, "11880.00"),
    @(160, "70.00"),
    @(161, "236.68"),
    @(162, "300.00"),
    @(163, "2589.00"),
    @(164, "7112.49"),
    @(165, "426.81"),
    @(166, "8672.13"),
    @(167, "51200.00"),
    @(168, "242.56"),
    @(169, "7398.99"),
    @(170, "1802.00"),
    @(171, "5910.00"),
    @(172, "14900.00"),
    @(173, "2700.00"),
    @(174, "1020.46"),
    @(175, "2587.64"),
    @(176, "4052.61"),
    @(177, "335.40"),
    @(178, "74709.92"),
    @(179, "125162.00"),
    @(180, "147876.00"),
    @(181, "700.00"),
    @(182, "32130.00"),
    @(183, "2637.16"),
    @(184, "1350.00"),
    @(185, "741.20"),
    @(186, "36000.00"),
    @(187, "1742.97"),
    @(188, "1192.23"),
    @(189, "8010.00"),
    @(190, "250512.74"),
    @(191, "262500.00"),
    @(192, "118000.00"),
    @(193, "185000.00"),
    @(194, "451142.00"),
    @(195, "194336.00"),
    @(196, "32500.00"),
    @(197, "27000.00"),
    @(198, "139500.00"),
    @(199, "300372.00"),
    @(200, "272800.00"),
    @(201, "38000.00"),
    @(202, "193700.00"),
    @(203, "92780.00"),
    @(204, "100000.00"),
    @(205, "14286789.66"),
    @(206, "70611.97"),
    @(207, "826292.78"),
    @(208, "30448.53"),
    @(209, "30120.00"),
    @(210, "3500.00"),
    @(211, "19840.00"),
    @(212, "3360.00"),
    @(213, "14973.00"),
    @(214, "1230.00"),
    @(215, "603.84"),
    @(216, "1200.01"),
    @(217, "1340.00"),
    @(218, "1325.00"),
    @(219, "755.00"),
    @(220, "680.00"),
    @(221, "30000.00")
)

$importeRange = $ws.Range("H2:H221")
$importeRange.NumberFormat = "@"

foreach ($row in $importeFixes) {
    $r = $row[0]
    $v = $row[1]
    $ws.Cells.Item($r, 8).Value = $v
}

$importeRange.ClearFormats()
